$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 and 31 swap places (Filecoin <-> Hedera) with refreshed price/volume data.
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05571"
$ws.Range("E30").Value = "  +1.80%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.125"
$ws.Range("E31").Value = "  -0.33%  "

# Refreshed Price (column D) and Volume(1h) (column E) values for every other coin row.
$updates = @(
    @{ Row = 2; D = "28.963.83"; E = "  -0.09%  " },
    @{ Row = 3; D = "1.827.29"; E = "  +0.06%  " },
    @{ Row = 4; D = "0.9946"; E = "  -0.17%  " },
    @{ Row = 5; D = "241.63"; E = "  -0.88%  " },
    @{ Row = 6; D = "0.6158"; E = "  -2.46%  " },
    @{ Row = 7; D = "0.9980"; E = "  -0.08%  " },
    @{ Row = 8; D = "0.07428"; E = "  -1.15%  " },
    @{ Row = 9; D = "0.2906"; E = "  -1.20%  " },
    @{ Row = 10; D = "22.93"; E = "  -0.59%  " },
    @{ Row = 11; D = "0.07632"; E = "  -0.84%  " },
    @{ Row = 12; D = "1.828.37"; E = "  -0.30%  " },
    @{ Row = 13; D = "4.976"; E = "  -0.31%  " },
    @{ Row = 14; D = "0.6705"; E = "  +0.27%  " },
    @{ Row = 15; D = "82.52"; E = "  -0.64%  " },
    @{ Row = 16; D = "0.000009179"; E = "  -4.67%  " },
    @{ Row = 17; D = "5.888"; E = "  -2.64%  " },
    @{ Row = 18; D = "28.969.18"; E = "  -0.22%  " },
    @{ Row = 19; D = "2.078.37"; E = "  -0.05%  " },
    @{ Row = 20; D = "239.64"; E = "  +5.97%  " },
    @{ Row = 21; D = "12.63"; E = "  +0.47%  " },
    @{ Row = 22; D = "0.9981"; E = "  +0.03%  " },
    @{ Row = 23; D = "7.186"; E = "  +0.57%  " },
    @{ Row = 24; D = "0.9954"; E = "  -0.25%  " },
    @{ Row = 25; D = "158.26"; E = "  -1.33%  " },
    @{ Row = 26; D = "0.1403"; E = "  -1.55%  " },
    @{ Row = 27; D = "8.462"; E = "  -0.48%  " },
    @{ Row = 28; D = "17.79"; E = "  -0.68%  " },
    @{ Row = 29; D = "1.496"; E = "  -0.31%  " },
    @{ Row = 32; D = "4.093"; E = "  +0.82%  " },
    @{ Row = 33; D = "1.202"; E = "  +0.13%  " },
    @{ Row = 34; D = "1.834"; E = "  -1.14%  " },
    @{ Row = 35; D = "0.7373"; E = "  -0.91%  " },
    @{ Row = 36; D = "1.136"; E = "  -0.11%  " },
    @{ Row = 37; D = "2.648"; E = "  +0.16%  " },
    @{ Row = 38; D = "2.762"; E = "  +0.38%  " },
    @{ Row = 39; D = "0.01779"; E = "  -0.04%  " },
    @{ Row = 40; D = "1.206.90"; E = "  -2.79%  " },
    @{ Row = 41; D = "6.411"; E = "  -3.57%  " },
    @{ Row = 42; D = "0.8949"; E = "  -0.52%  " },
    @{ Row = 43; D = "0.9961"; E = "  -0.27%  " },
    @{ Row = 44; D = "101.03"; E = "  -0.26%  " },
    @{ Row = 45; D = "1.976.55"; E = "  -0.09%  " },
    @{ Row = 46; D = "65.16"; E = "  -0.15%  " },
    @{ Row = 47; D = "0.5067"; E = "  -0.38%  " },
    @{ Row = 48; D = "0.00000000118"; E = "  -5.05%  " },
    @{ Row = 49; D = "0.4045"; E = "  -0.20%  " },
    @{ Row = 50; D = "9.121"; E = "  +2.04%  " },
    @{ Row = 51; D = "0.05799"; E = "  +0.26%  " }
)

foreach ($u in $updates) {
    # Force column D to stay text so numeric-looking values (e.g. "0.9946")
    # are not auto-converted to numbers by Excel, matching the original
    # inline-string (text) cell content.
    $ws.Range("D" + $u.Row).NumberFormat = "@"
    $ws.Range("D" + $u.Row).Value = $u.D
    $ws.Range("E" + $u.Row).Value = $u.E
}
